# Applies the commit "Atualizado por script em 01-12-2023 20:45":
#  1) Swap match details (cols F:V) between rows 23 and 24
#  2) Swap match details (cols F:V) between rows 74 and 75
#  3) Append two new match rows (101 and 102)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Swap rows 23 and 24 (columns F..V) ---
$row23 = $ws.Range("F23:V23").Value2
$row24 = $ws.Range("F24:V24").Value2
$ws.Range("F23:V23").Value2 = $row24
$ws.Range("F24:V24").Value2 = $row23

# --- 2) Swap rows 74 and 75 (columns F..V) ---
$row74 = $ws.Range("F74:V74").Value2
$row75 = $ws.Range("F75:V75").Value2
$ws.Range("F74:V74").Value2 = $row75
$ws.Range("F75:V75").Value2 = $row74

# --- 3) Append new rows 101 and 102, copying formatting from row 100 ---
$ws.Range("A100:V100").Copy()
$ws.Range("A101:V102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 101: Vilaverdense 1 - 0 Leiria
$ws.Range("A101").Value = 100
$ws.Range("B101").Value = "portugal"
$ws.Range("C101").Value = "liga-portugal-2"
$ws.Range("D101").Value = "2023-2024"
$ws.Range("E101").Value2 = 45261.6875
$ws.Range("F101").Value = "Vilaverdense"
$ws.Range("G101").Value = 1
$ws.Range("H101").Value = "Leiria"
$ws.Range("I101").Value = 0
$ws.Range("J101").Value = 3.9
$ws.Range("K101").Value = "24/11/2023 16:43"
$ws.Range("L101").Value = 4.05
$ws.Range("M101").Value = "01/12/2023 14:37"
$ws.Range("N101").Value = 3.76
$ws.Range("O101").Value = "24/11/2023 16:43"
$ws.Range("P101").Value = 3.8
$ws.Range("Q101").Value = "01/12/2023 15:02"
$ws.Range("R101").Value = 1.85
$ws.Range("S101").Value = "24/11/2023 16:43"
$ws.Range("T101").Value = 1.9
$ws.Range("U101").Value = "01/12/2023 14:37"
$ws.Range("V101").Value = "https://www.betexplorer.com/football/portugal/liga-portugal-2/vilaverdense-fc-leiria/CIbQlxdr/"

# Row 102: Academico Viseu 3 - 1 Os Belenenses
$ws.Range("A102").Value = 101
$ws.Range("B102").Value = "portugal"
$ws.Range("C102").Value = "liga-portugal-2"
$ws.Range("D102").Value = "2023-2024"
$ws.Range("E102").Value2 = 45261.79166666666
$ws.Range("F102").Value = "Academico Viseu"
$ws.Range("G102").Value = 3
$ws.Range("H102").Value = "Os Belenenses"
$ws.Range("I102").Value = 1
$ws.Range("J102").Value = 1.56
$ws.Range("K102").Value = "25/11/2023 12:12"
$ws.Range("L102").Value = 1.65
$ws.Range("M102").Value = "01/12/2023 18:59"
$ws.Range("N102").Value = 4.12
$ws.Range("O102").Value = "25/11/2023 12:12"
$ws.Range("P102").Value = 3.84
$ws.Range("Q102").Value = "01/12/2023 18:59"
$ws.Range("R102").Value = 5.37
$ws.Range("S102").Value = "25/11/2023 12:12"
$ws.Range("T102").Value = 5.82
$ws.Range("U102").Value = "01/12/2023 18:59"
$ws.Range("V102").Value = "https://www.betexplorer.com/football/portugal/liga-portugal-2/academico-viseu-cf-os-belenenses/j9cUmdBl/"
